$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Saps": insert a new first column "SapID" holding the original
# MongoDB ObjectId string for every sap record.
# ---------------------------------------------------------------------------
$saps = $wb.Worksheets.Item("Saps")
$saps.Columns.Item(1).Insert()
$saps.Columns.Item(1).ColumnWidth = $saps.Columns.Item(2).ColumnWidth

$saps.Cells.Item(1, 1).Value2 = "SapID"

$sapIds = @(
    '"60776b61ccab402de07f4c81"',
    '"60776e12ccab402de07f4c82"',
    '"607787dc0eadad1a7868fe6c"',
    '"607789ac06902a1e3881b27a"',
    '"60778a8bd5577745f4085a06"',
    '"60778b14e0edf5397cb584f6"',
    '"60778bec5678c12dfc43ae5c"',
    '"60778caa99ca002c28453020"',
    '"60778d3663098e3e98890744"',
    '"60778d889f9f3a3c6447bc04"',
    '"6077974db6508a1cd0e0f073"'
)

for ($i = 0; $i -lt $sapIds.Length; $i++) {
    $row = $i + 2
    $saps.Cells.Item($row, 1).Value2 = $sapIds[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Syrups": insert two new leading columns - "User" and "Syrup ID" -
# holding the owning user and the original MongoDB ObjectId string for
# every syrup record.
# ---------------------------------------------------------------------------
$syrups = $wb.Worksheets.Item("Syrups")
$syrups.Range("A1:B1").EntireColumn.Insert()
$syrups.Columns.Item(1).ColumnWidth = $syrups.Columns.Item(3).ColumnWidth
$syrups.Columns.Item(2).ColumnWidth = $syrups.Columns.Item(3).ColumnWidth

$syrups.Cells.Item(1, 1).Value2 = "User"
$syrups.Cells.Item(1, 2).Value2 = "Syrup ID"

$syrupIds = @(
    '"605e40cd60e99e32004bb4dc"',
    '"6064f27d5320e055844a8c05"',
    '"60679b88c4170d31d476a7bb"',
    '"6070b6aa7d0400083c267385"'
)

for ($i = 0; $i -lt $syrupIds.Length; $i++) {
    $row = $i + 2
    $syrups.Cells.Item($row, 1).Value2 = "tang"
    $syrups.Cells.Item($row, 2).Value2 = $syrupIds[$i]
}
